$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '20.563.07'
$ws.Range('D3').Value = '1.470.29'
$ws.Range('E3').Value = '  +2.03%  '
$ws.Range('E4').Value = '  +0.26%  '
$ws.Range('D5').Value = "'0.9588"
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +4.09%  '
$ws.Range('D6').Value = "'276.68"
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.69%  '
$ws.Range('D7').Value = "'0.3558"
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -2.16%  '
$ws.Range('D8').Value = "'0.3058"
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.35%  '
$ws.Range('D9').Value = "'1.082"
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +6.73%  '
$ws.Range('D10').Value = "'39.16"
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +1.14%  '
$ws.Range('D11').Value = "'0.06611"
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +2.43%  '
$ws.Range('E12').Value = '  +0.41%  '
$ws.Range('D13').Value = "'5.444"
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +2.78%  '
$ws.Range('D14').Value = "'18.03"
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +4.15%  '
$ws.Range('E15').Value = '  +2.67%  '
$ws.Range('D16').Value = "'0.9602"
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +2.27%  '
$ws.Range('E17').Value = '  +1.34%  '
$ws.Range('D18').Value = '1.468.17'
$ws.Range('E18').Value = '  +1.81%  '
$ws.Range('D19').Value = "'0.05929"
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +5.65%  '
$ws.Range('D20').Value = "'68.68"
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +1.78%  '
$ws.Range('D21').Value = "'5.460"
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +2.96%  '
$ws.Range('D22').Value = "'14.42"
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +2.08%  '
$ws.Range('E23').Value = '  +4.29%  '
$ws.Range('D24').Value = "'2.275"
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +1.49%  '
$ws.Range('D25').Value = '20.572.18'
$ws.Range('E25').Value = '  +1.50%  '
$ws.Range('D26').Value = "'145.73"
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +5.17%  '
$ws.Range('D27').Value = "'2.078"
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +2.38%  '
$ws.Range('D28').Value = "'17.05"
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +1.05%  '
$ws.Range('D29').Value = '1.631.60'
$ws.Range('E29').Value = '  +2.29%  '
$ws.Range('D30').Value = "'114.16"
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +3.78%  '
$ws.Range('D31').Value = "'4.000"
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.38%  '
$ws.Range('D32').Value = "'4.901"
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +2.59%  '
$ws.Range('D33').Value = "'0.07923"
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +3.65%  '
$ws.Range('D34').Value = "'0.7872"
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +2.35%  '
$ws.Range('D35').Value = "'1.211"
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +7.00%  '
$ws.Range('D36').Value = "'1.452"
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.07%  '
$ws.Range('D37').Value = "'0.05661"
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -1.03%  '
$ws.Range('D38').Value = "'4.712"
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +2.32%  '
$ws.Range('D39').Value = "'0.9605"
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +3.34%  '
$ws.Range('E40').Value = '  +2.55%  '
$ws.Range('E41').Value = '  +1.48%  '
$ws.Range('D42').Value = "'0.1840"
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +0.72%  '
$ws.Range('D43').Value = "'7.265"
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +4.98%  '
$ws.Range('D44').Value = "'3.510"
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.00%  '
$ws.Range('D45').Value = "'0.5207"
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.89%  '
$ws.Range('D46').Value = "'12.06"
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +3.06%  '
$ws.Range('D47').Value = "'119.51"
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +4.19%  '
$ws.Range('D48').Value = "'0.5153"
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +2.52%  '
$ws.Range('E49').Value = '  +4.62%  '
$ws.Range('D50').Value = "'0.06416"
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.00%  '
$ws.Range('D51').Value = "'0.9938"
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.86%  '
